$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$null = $tr.InsertAfter("`rWrite your learning notes (added to ")
$null = $tr.InsertAfter("resume later)")
